$wb = $excel.ActiveWorkbook

# This script re-applies refreshed Universalis market-price snapshots
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ)
# and the dependent LevePrice/LeveProfit columns for each Sheet tab,
# mirroring the scheduled-runner price refresh described in the commit.

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1907.8276
$ws.Range("J17").Value = 1907.8276
$ws.Range("L17").Value = 5723.4828
$ws.Range("N17").Value = -6059.4828
# Row 55
$ws.Range("H55").Value = 300
$ws.Range("J55").Value = 300
$ws.Range("L55").Value = 300
$ws.Range("N55").Value = -728
# Row 64
$ws.Range("H64").Value = 4444.024
$ws.Range("I64").Value = 3741.4348
$ws.Range("K64").Value = 3741.4348
$ws.Range("M64").Value = -3493.4348
# Row 67
$ws.Range("H67").Value = 4444.024
$ws.Range("I67").Value = 3741.4348
$ws.Range("K67").Value = 3741.4348
$ws.Range("M67").Value = -2883.4348
# Row 112
$ws.Range("H112").Value = 979.3333
$ws.Range("I112").Value = 972.6667
$ws.Range("J112").Value = 981
$ws.Range("K112").Value = 2918.0001
$ws.Range("L112").Value = 2943
$ws.Range("M112").Value = -1810.0001
$ws.Range("N112").Value = -5159
# Row 113
$ws.Range("H113").Value = 3664.4707
$ws.Range("J113").Value = 3736.3635
$ws.Range("L113").Value = 3736.3635
$ws.Range("N113").Value = -10244.3635
# Row 135
$ws.Range("H135").Value = 1920.8636
$ws.Range("I135").Value = 1653.4445
$ws.Range("K135").Value = 14881.0005
$ws.Range("M135").Value = -12346.0005
# Row 138
$ws.Range("H138").Value = 2246.3088
$ws.Range("I138").Value = 1534.6154
$ws.Range("J138").Value = 2414.5273
$ws.Range("K138").Value = 4603.8462
$ws.Range("L138").Value = 7243.581900000001
$ws.Range("M138").Value = 536.1538
$ws.Range("N138").Value = -17523.5819
# Row 141
$ws.Range("H141").Value = 4358.75
$ws.Range("I141").Value = 3645
$ws.Range("J141").Value = 6500
$ws.Range("K141").Value = 10935
$ws.Range("L141").Value = 19500
$ws.Range("M141").Value = -5755
$ws.Range("N141").Value = -29860

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6953.8364
$ws.Range("I32").Value = 3400.75
$ws.Range("K32").Value = 3400.75
$ws.Range("M32").Value = -3113.75
# Row 45
$ws.Range("H45").Value = 6251996
$ws.Range("I45").Value = 1993.3334
$ws.Range("J45").Value = 15627000
$ws.Range("K45").Value = 1993.3334
$ws.Range("L45").Value = 15627000
$ws.Range("M45").Value = -1616.3334
$ws.Range("N45").Value = -15627754
# Row 61
$ws.Range("H61").Value = 61621.234
$ws.Range("I61").Value = 2736
$ws.Range("J61").Value = 252998.25
$ws.Range("K61").Value = 2736
$ws.Range("L61").Value = 252998.25
$ws.Range("M61").Value = -2524
$ws.Range("N61").Value = -253422.25
# Row 74
$ws.Range("H74").Value = 65727.25
$ws.Range("I74").Value = 144912.58
$ws.Range("K74").Value = 144912.58
$ws.Range("M74").Value = -144038.58
# Row 77
$ws.Range("H77").Value = 65727.25
$ws.Range("I77").Value = 144912.58
$ws.Range("K77").Value = 724562.8999999999
$ws.Range("M77").Value = -720194.8999999999
# Row 107
$ws.Range("H107").Value = 46664.668
$ws.Range("J107").Value = 46664.668
$ws.Range("L107").Value = 46664.668
$ws.Range("N107").Value = -54344.668
# Row 132
$ws.Range("H132").Value = 2901.1428
$ws.Range("I132").Value = 2457.3333
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 7371.999899999999
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -4841.999899999999
$ws.Range("N132").Value = -16160
# Row 136
$ws.Range("H136").Value = 61621.234
$ws.Range("I136").Value = 2736
$ws.Range("J136").Value = 252998.25
$ws.Range("K136").Value = 8208
$ws.Range("L136").Value = 758994.75
$ws.Range("M136").Value = -5658
$ws.Range("N136").Value = -764094.75

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 3064.6
$ws.Range("I54").Value = 1738.4445
$ws.Range("K54").Value = 1738.4445
$ws.Range("M54").Value = -1254.4445
# Row 86
$ws.Range("H86").Value = 7130.3335
$ws.Range("I86").Value = 4172.091
$ws.Range("J86").Value = 11779
$ws.Range("K86").Value = 4172.091
$ws.Range("L86").Value = 11779
$ws.Range("M86").Value = -3049.091
$ws.Range("N86").Value = -14025
# Row 89
$ws.Range("H89").Value = 7130.3335
$ws.Range("I89").Value = 4172.091
$ws.Range("J89").Value = 11779
$ws.Range("K89").Value = 20860.455
$ws.Range("L89").Value = 58895
$ws.Range("M89").Value = -15244.455
$ws.Range("N89").Value = -70127
# Row 134
$ws.Range("H134").Value = 2515.762
$ws.Range("I134").Value = 1046.1666
$ws.Range("J134").Value = 11333.333
$ws.Range("K134").Value = 3138.4998
$ws.Range("L134").Value = 33999.999
$ws.Range("M134").Value = -603.4998000000001
$ws.Range("N134").Value = -39069.999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2039.5245
$ws.Range("I31").Value = 1434.3556
$ws.Range("K31").Value = 1434.3556
$ws.Range("M31").Value = -1139.3556
# Row 34
$ws.Range("H34").Value = 2039.5245
$ws.Range("I34").Value = 1434.3556
$ws.Range("K34").Value = 1434.3556
$ws.Range("M34").Value = -1232.3556
# Row 132
$ws.Range("H132").Value = 2382.7144
$ws.Range("I132").Value = 2205.4443
$ws.Range("J132").Value = 2701.8
$ws.Range("K132").Value = 6616.3329
$ws.Range("L132").Value = 8105.400000000001
$ws.Range("M132").Value = -4086.3329
$ws.Range("N132").Value = -13165.4
# Row 134
$ws.Range("H134").Value = 33544.934
$ws.Range("I134").Value = 950
$ws.Range("J134").Value = 85153.586
$ws.Range("K134").Value = 2850
$ws.Range("L134").Value = 255460.758
$ws.Range("M134").Value = -315
$ws.Range("N134").Value = -260530.758

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 301.2
$ws.Range("I10").Value = 333
$ws.Range("K10").Value = 999
$ws.Range("M10").Value = -860
# Row 36
$ws.Range("H36").Value = 3251
$ws.Range("I36").Value = 3251
$ws.Range("K36").Value = 9753
$ws.Range("M36").Value = -9584
# Row 136
$ws.Range("H136").Value = 2215.875
$ws.Range("I136").Value = 1888.9231
$ws.Range("K136").Value = 5666.7693
$ws.Range("M136").Value = -566.7692999999999
# Row 138
$ws.Range("H138").Value = 6417.6772
$ws.Range("J138").Value = 6235.25
$ws.Range("L138").Value = 18705.75
$ws.Range("N138").Value = -28985.75

$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 2166.6667
$ws.Range("I36").Value = 750
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 750
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -265
$ws.Range("N36").Value = -5970
# Row 70
$ws.Range("H70").Value = 155715.17
$ws.Range("I70").Value = 7425.6665
$ws.Range("K70").Value = 7425.6665
$ws.Range("M70").Value = -7155.6665
# Row 73
$ws.Range("H73").Value = 155715.17
$ws.Range("I73").Value = 7425.6665
$ws.Range("K73").Value = 7425.6665
$ws.Range("M73").Value = -6489.6665
# Row 102
$ws.Range("H102").Value = 1617.875
$ws.Range("I102").Value = 1518.7894
$ws.Range("K102").Value = 1518.7894
$ws.Range("M102").Value = 103.2106000000001
# Row 107
$ws.Range("H107").Value = 917.0476
$ws.Range("I107").Value = 710.7143
$ws.Range("J107").Value = 1020.2143
$ws.Range("K107").Value = 710.7143
$ws.Range("L107").Value = 1020.2143
$ws.Range("M107").Value = 1209.2857
$ws.Range("N107").Value = -4860.2143
# Row 132
$ws.Range("H132").Value = 4162.619
$ws.Range("I132").Value = 2708.6428
$ws.Range("J132").Value = 7070.5713
$ws.Range("K132").Value = 8125.928400000001
$ws.Range("L132").Value = 21211.7139
$ws.Range("M132").Value = -5595.928400000001
$ws.Range("N132").Value = -26271.7139
# Row 135
$ws.Range("H135").Value = 56340.625
$ws.Range("J135").Value = 56340.625
$ws.Range("L135").Value = 56340.625
$ws.Range("N135").Value = -66480.625

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 901.0833
$ws.Range("I16").Value = 910.0909
$ws.Range("K16").Value = 910.0909
$ws.Range("M16").Value = -740.0909
# Row 22
$ws.Range("H22").Value = 8591.923000000001
$ws.Range("I22").Value = 1462.5
$ws.Range("J22").Value = 11760.556
$ws.Range("K22").Value = 1462.5
$ws.Range("L22").Value = 11760.556
$ws.Range("M22").Value = -1167.5
$ws.Range("N22").Value = -12350.556
# Row 27
$ws.Range("H27").Value = 8591.923000000001
$ws.Range("I27").Value = 1462.5
$ws.Range("J27").Value = 11760.556
$ws.Range("K27").Value = 1462.5
$ws.Range("L27").Value = 11760.556
$ws.Range("M27").Value = -1355.5
$ws.Range("N27").Value = -11974.556
# Row 46
$ws.Range("H46").Value = 22259
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376
# Row 122
$ws.Range("H122").Value = 75004230
$ws.Range("I122").Value = 71433140
$ws.Range("J122").Value = 100001900
$ws.Range("K122").Value = 214299420
$ws.Range("L122").Value = 300005700
$ws.Range("M122").Value = -214296970
$ws.Range("N122").Value = -300010600
# Row 132
$ws.Range("H132").Value = 2865.5
$ws.Range("I132").Value = 2603.6924
$ws.Range("K132").Value = 7811.0772
$ws.Range("M132").Value = -5281.0772

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1260.5555
$ws.Range("I113").Value = 906.7692
$ws.Range("K113").Value = 2720.3076
$ws.Range("M113").Value = -550.3076000000001
# Row 122
$ws.Range("H122").Value = 1917.2941
$ws.Range("I122").Value = 1237.1111
$ws.Range("K122").Value = 3711.3333
$ws.Range("M122").Value = -1261.3333
# Row 132
$ws.Range("H132").Value = 2078.182
$ws.Range("I132").Value = 1743.7916
$ws.Range("J132").Value = 2969.889
$ws.Range("K132").Value = 5231.3748
$ws.Range("L132").Value = 8909.667000000001
$ws.Range("M132").Value = -2701.3748
$ws.Range("N132").Value = -13969.667
# Row 136
$ws.Range("H136").Value = 2054.8462
$ws.Range("I136").Value = 1833
$ws.Range("J136").Value = 2193.5
$ws.Range("K136").Value = 5499
$ws.Range("L136").Value = 6580.5
$ws.Range("M136").Value = -2949
$ws.Range("N136").Value = -11680.5
